$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("NAQUICHE MECHATO ENMA", 131),
    @("NAQUICHE SILVA MARIA LIZETH", 103),
    @("CUBAS GARCIA ROSA ANITA", 76),
    @("MANOSALVA RUIZ SANDRA KAROLINE", 73),
    @("PACHECO ALISON", 72),
    @("CORAS QUISPE JORGE AMERICO", 71),
    @("CASTREJON TELLO GRECIA", 69),
    @("BECERRA ASMAT CAROL STEFANY", 61),
    @("SAUCEDO CABRERA CARLOS ALEXANDER", 58)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}
